$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.592.91"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "2.732.28"
$ws.Range("E3").Value = "  +4.20%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'526.24"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").Value = "'145.63"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.577"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "2.731.59"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").Value = "'6.78"
$ws.Range("E10").Value = "  +7.21%  "
$ws.Range("D11").Value = "'0.106"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("E13").Value = "  +3.27%  "
$ws.Range("D14").Value = "3.181.23"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "60.584.02"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.921.31"
$ws.Range("E16").Value = "  +10.81%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'21.29"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "'344.70"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "'10.63"
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("E22").Value = "  +4.88%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D24").Value = "'63.36"
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("D25").Value = "'0.421"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").Value = "'0.169"
$ws.Range("E26").Value = "  +3.04%  "
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "0.0₃0820"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").Value = "'7.28"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").Value = "'6.84"
$ws.Range("E30").Value = "  +9.57%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").Value = "'149.53"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").Value = "'4.27"
$ws.Range("E35").Value = "  +7.52%  "
$ws.Range("E36").Value = "  +8.12%  "
$ws.Range("D37").Value = "'0.938"
$ws.Range("E37").Value = "  -4.04%  "
$ws.Range("D38").Value = "'0.878"
$ws.Range("E38").Value = "  +4.56%  "
$ws.Range("D39").Value = "'1.53"
$ws.Range("E39").Value = "  +7.88%  "
$ws.Range("D40").Value = "'37.08"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "'282.26"
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("D43").Value = "'20.16"
$ws.Range("E43").Value = "  +3.72%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.612"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "2.146.12"
$ws.Range("E46").Value = "  +8.05%  "
$ws.Range("D47").Value = "'0.0986"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "'4.91"
$ws.Range("E48").Value = "  +6.04%  "
$ws.Range("D49").Value = "'0.0539"
$ws.Range("E49").Value = "  +3.27%  "
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").Value = "'0.0232"
$ws.Range("E51").Value = "  +1.52%  "
